$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Session 3 (Divide & Conquer) mark for the student
$ws.Range("D4").Value = 5

# Teacher's feedback comment for the D&C session
$ws.Range("D5").Value = "The D&C version is not working as expected. Please, check the video of the last seminar (the implementation should be very similar to Mergesort). Second PDF is missing"

# Update the active selection to reflect the cell that was just edited
$ws.Range("D5:D12").Select()
